# Apply odds updates to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    "N2"  = 4.3
    "J3"  = 3.55
    "Q3"  = 1.83
    "AD3" = 24
    "AF3" = 18
    "AH3" = 25
    "J4"  = 4.5
    "F5"  = 2.1
    "G5"  = 2.38
    "H5"  = 3.4
    "I5"  = 4.5
    "J5"  = 3.2
    "N5"  = 3.2
    "O5"  = 1.34
    "P5"  = 1.75
    "Q5"  = 2.2
    "R5"  = 1.2
    "S5"  = 3.55
    "T5"  = 1.01
    "U5"  = 1.97
    "W5"  = 1.72
    "X5"  = 1000
    "Y5"  = 19
    "Z5"  = 40
    "AA5" = 100
    "AB5" = 13
    "AC5" = 11
    "AD5" = 23
    "AE5" = 70
    "AF5" = 20
    "AG5" = 16
    "AH5" = 28
    "AI5" = 90
    "AJ5" = 44
    "AK5" = 38
    "AL5" = 980
    "AN5" = 30
    "F6"  = 1.86
    "I6"  = 5.4
    "U6"  = 1.75
    "V6"  = 1.23
    "F7"  = 1.74
    "G7"  = 2.04
    "I7"  = 6.6
    "J7"  = 3.2
    "Q7"  = 1.77
    "U7"  = 1.98
    "W7"  = 1.98
    "AC7" = 980
    "G8"  = 2.52
    "H8"  = 3.25
    "K8"  = 3.6
    "L8"  = 1.45
    "N8"  = 3.15
    "Q8"  = 1.96
    "R8"  = 1.28
    "U8"  = 1.99
    "W8"  = 1.66
    "X8"  = 14.5
    "Y8"  = 14.5
    "Z8"  = 29
    "AA8" = 80
    "AB8" = 11
    "AC8" = 9.199999999999999
    "AD8" = 18
    "AE8" = 55
    "AF8" = 18
    "AG8" = 14
    "AH8" = 23
    "AI8" = 70
    "AJ8" = 42
    "AK8" = 36
    "AL8" = 55
    "AN8" = 30
    "G9"  = 2.18
    "H9"  = 4.1
    "I9"  = 5.5
    "N9"  = 1.01
    "O9"  = 1.36
    "P9"  = 1.76
    "Q9"  = 1.94
    "W9"  = 1.85
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
